# Fix the casing of two registration-number entries that were typed in
# lower-case, and restore the window's scroll/selection state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct "23bce2327" -> "23BCE2327" (row 106) first so that the new
# shared-string entries are appended in the same order Excel produced them.
$ws.Range("A106").Value = "23BCE2327"

# Correct "23mis0079" -> "23MIS0079" (row 48)
$ws.Range("A48").Value = "23MIS0079"

# Restore the sheet view: scrolled so row 67 is at the top, with D95 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("D95").Select()
